$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '27.216.05'
Set-TextValue $ws.Range('E2') '  +0.51%  '

Set-TextValue $ws.Range('D3') '1.894.47'
Set-TextValue $ws.Range('E3') '  -0.01%  '

Set-TextValue $ws.Range('D4') '1.003'
Set-TextValue $ws.Range('E4') '  +0.15%  '

Set-TextValue $ws.Range('D5') '307.56'
Set-TextValue $ws.Range('E5') '  +0.15%  '

Set-TextValue $ws.Range('E6') '  +0.22%  '

Set-TextValue $ws.Range('D7') '0.5176'
Set-TextValue $ws.Range('E7') '  -0.13%  '

Set-TextValue $ws.Range('D8') '0.3755'
Set-TextValue $ws.Range('E8') '  -0.12%  '

Set-TextValue $ws.Range('D9') '0.07272'
Set-TextValue $ws.Range('E9') '  +0.67%  '

Set-TextValue $ws.Range('D10') '21.16'
Set-TextValue $ws.Range('E10') '  -0.05%  '

Set-TextValue $ws.Range('D11') '0.8995'
Set-TextValue $ws.Range('E11') '  +0.62%  '

Set-TextValue $ws.Range('D12') '0.08128'
Set-TextValue $ws.Range('E12') '  +6.08%  '

Set-TextValue $ws.Range('D13') '96.27'
Set-TextValue $ws.Range('E13') '  +2.15%  '

Set-TextValue $ws.Range('D14') '1.896.50'
Set-TextValue $ws.Range('E14') '  +0.06%  '

Set-TextValue $ws.Range('D15') '5.280'
Set-TextValue $ws.Range('E15') '  +1.02%  '

Set-TextValue $ws.Range('D16') '1.003'
Set-TextValue $ws.Range('E16') '  +0.22%  '

Set-TextValue $ws.Range('D17') '0.000008591'
Set-TextValue $ws.Range('E17') '  +0.91%  '

Set-TextValue $ws.Range('D18') '14.53'
Set-TextValue $ws.Range('E18') '  +0.32%  '

Set-TextValue $ws.Range('E19') '  +0.16%  '

Set-TextValue $ws.Range('D20') '27.239.85'
Set-TextValue $ws.Range('E20') '  +0.44%  '

Set-TextValue $ws.Range('D21') '5.079'
Set-TextValue $ws.Range('E21') '  +0.34%  '

Set-TextValue $ws.Range('D22') '10.67'
Set-TextValue $ws.Range('E22') '  +0.76%  '

Set-TextValue $ws.Range('D23') '6.388'
Set-TextValue $ws.Range('E23') '  -0.37%  '

Set-TextValue $ws.Range('D24') '2.297'
Set-TextValue $ws.Range('E24') '  +0.30%  '

Set-TextValue $ws.Range('D25') '147.06'
Set-TextValue $ws.Range('E25') '  +0.40%  '

Set-TextValue $ws.Range('B26') 'Toncoin'
Set-TextValue $ws.Range('C26') 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range('D26') '1.745'
Set-TextValue $ws.Range('E26') '  +0.71%  '

Set-TextValue $ws.Range('B27') 'EthereumClassic'
Set-TextValue $ws.Range('C27') 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range('D27') '18.20'
Set-TextValue $ws.Range('E27') '  +0.82%  '

Set-TextValue $ws.Range('D28') '115.09'
Set-TextValue $ws.Range('E28') '  +0.52%  '

Set-TextValue $ws.Range('B29') 'Filecoin'
Set-TextValue $ws.Range('C29') 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range('D29') '4.947'
Set-TextValue $ws.Range('E29') '  -0.55%  '

Set-TextValue $ws.Range('B30') 'InternetComputer(DFINITY)'
Set-TextValue $ws.Range('C30') 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range('D30') '4.822'
Set-TextValue $ws.Range('E30') '  +0.89%  '

Set-TextValue $ws.Range('D31') '0.09223'
Set-TextValue $ws.Range('E31') '  +0.27%  '

Set-TextValue $ws.Range('B32') 'Hedera'
Set-TextValue $ws.Range('C32') 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D32') '0.05027'
Set-TextValue $ws.Range('E32') '  -0.38%  '

Set-TextValue $ws.Range('B33') 'ImmutableX'
Set-TextValue $ws.Range('C33') 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range('D33') '0.7918'
Set-TextValue $ws.Range('E33') '  +2.11%  '

Set-TextValue $ws.Range('D34') '1.216'
Set-TextValue $ws.Range('E34') '  -1.72%  '

Set-TextValue $ws.Range('E35') '  +4.96%  '

Set-TextValue $ws.Range('D36') '2.966'
Set-TextValue $ws.Range('E36') '  -0.33%  '

Set-TextValue $ws.Range('D37') '2.589'
Set-TextValue $ws.Range('E37') '  +0.00%  '

Set-TextValue $ws.Range('D38') '0.5648'
Set-TextValue $ws.Range('E38') '  +0.68%  '

Set-TextValue $ws.Range('D39') '0.01982'
Set-TextValue $ws.Range('E39') '  -0.27%  '

Set-TextValue $ws.Range('D40') '1.073'
Set-TextValue $ws.Range('E40') '  -0.16%  '

Set-TextValue $ws.Range('D41') '8.938'
Set-TextValue $ws.Range('E41') '  -0.45%  '

Set-TextValue $ws.Range('D42') '6.543'
Set-TextValue $ws.Range('E42') '  -1.17%  '

Set-TextValue $ws.Range('D43') '115.28'
Set-TextValue $ws.Range('E43') '  -3.49%  '

Set-TextValue $ws.Range('D44') '0.1511'
Set-TextValue $ws.Range('E44') '  -0.37%  '

Set-TextValue $ws.Range('D45') '0.4849'
Set-TextValue $ws.Range('E45') '  +0.50%  '

Set-TextValue $ws.Range('E46') '  +0.22%  '

Set-TextValue $ws.Range('D47') '10.01'
Set-TextValue $ws.Range('E47') '  -1.29%  '

Set-TextValue $ws.Range('D48') '1.617'
Set-TextValue $ws.Range('E48') '  +1.44%  '

Set-TextValue $ws.Range('D49') '38.15'
Set-TextValue $ws.Range('E49') '  +1.84%  '

Set-TextValue $ws.Range('D50') '63.30'
Set-TextValue $ws.Range('E50') '  -0.93%  '

Set-TextValue $ws.Range('D51') '0.05945'
Set-TextValue $ws.Range('E51') '  +0.38%  '
